$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns keep their text formatting
# so numeric-looking strings such as "1.11" are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '59.302.72'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '2.519.65'
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '534.88'
$ws.Range("E5").Value = '  -1.34%  '
$ws.Range("D6").Value = '139.30'
$ws.Range("E6").Value = '  -3.78%  '
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("D8").Value = '0.565'
$ws.Range("E8").Value = '  -1.76%  '
$ws.Range("D9").Value = '2.524.75'
$ws.Range("E9").Value = '  -1.02%  '
$ws.Range("E10").Value = '  +0.07%  '
$ws.Range("E11").Value = '  +1.26%  '
$ws.Range("D12").Value = '5.45'
$ws.Range("E12").Value = '  -2.80%  '
$ws.Range("E13").Value = '  -0.10%  '
$ws.Range("D14").Value = '2.969.00'
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("D15").Value = '23.53'
$ws.Range("E15").Value = '  -1.30%  '
$ws.Range("D16").Value = '59.241.41'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").Value = '2.521.86'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("E19").Value = '  -1.32%  '
$ws.Range("D20").Value = '4.31'
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").Value = '325.16'
$ws.Range("E21").Value = '  -0.44%  '
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("D23").Value = '5.81'
$ws.Range("E23").Value = '  -0.65%  '
$ws.Range("D24").Value = '63.72'
$ws.Range("E24").Value = '  +2.48%  '
$ws.Range("D25").Value = '0.428'
$ws.Range("E25").Value = '  -2.38%  '
$ws.Range("E26").Value = '  +1.11%  '
$ws.Range("E27").Value = '  +0.95%  '
$ws.Range("D28").Value = '7.84'
$ws.Range("E28").Value = '  -2.11%  '
$ws.Range("D29").Value = '6.93'
$ws.Range("E29").Value = '  +0.95%  '
$ws.Range("D30").Value = '0.0₃0778'
$ws.Range("E30").Value = '  -0.93%  '
$ws.Range("D32").Value = '164.79'
$ws.Range("E32").Value = '  +4.87%  '
$ws.Range("E33").Value = '  -2.00%  '
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.17%  '
$ws.Range("D35").Value = '1.11'
$ws.Range("E35").Value = '  -9.78%  '
$ws.Range("D36").Value = '18.53'
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("D37").Value = '4.27'
$ws.Range("E37").Value = '  -3.02%  '
$ws.Range("D38").Value = '1.58'
$ws.Range("E38").Value = '  -2.17%  '
$ws.Range("D39").Value = '36.97'
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("D40").Value = '3.69'
$ws.Range("E40").Value = '  -0.65%  '
$ws.Range("D41").Value = '0.815'
$ws.Range("E41").Value = '  -1.95%  '
$ws.Range("D42").Value = '5.26'
$ws.Range("E42").Value = '  -7.51%  '
$ws.Range("D43").Value = '279.77'
$ws.Range("E43").Value = '  -6.65%  '
$ws.Range("E44").Value = '  +0.50%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '0.598'
$ws.Range("E45").Value = '  -1.14%  '
$ws.Range("B46").Value = 'WhiteBITCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D46").Value = '10.84'
$ws.Range("E46").Value = '  +0.41%  '
$ws.Range("D47").Value = '0.0935'
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("D48").Value = '123.16'
$ws.Range("E48").Value = '  +0.14%  '
$ws.Range("D49").Value = '0.0514'
$ws.Range("E49").Value = '  -0.71%  '
$ws.Range("E50").Value = '  -1.77%  '
$ws.Range("D51").Value = '17.78'
$ws.Range("E51").Value = '  -3.25%  '
